$wb = $excel.ActiveWorkbook

# --- Sheet "Ingreso": append 10 new aporte rows dated 2023-07-02 (serial 45109) ---
$wsIngreso = $wb.Worksheets.Item("Ingreso")

$aportes = @(
    @("Wilkin", 100),
    @("Carlos", 100),
    @("Anuel", 100),
    @("Frandy", 500),
    @("Chamo", 100),
    @("Invitados", 100),
    @("Kukito", 50),
    @("Orlando", 1000),
    @("Joel", 300),
    @("Mac Daniel", 900)
)

$startRow = 464
for ($i = 0; $i -lt $aportes.Count; $i++) {
    $row = $startRow + $i
    $wsIngreso.Cells.Item($row, 1).Value = 45109
    $wsIngreso.Cells.Item($row, 2).Value = $aportes[$i][0]
    $wsIngreso.Cells.Item($row, 3).Value = $aportes[$i][1]
    $wsIngreso.Cells.Item($row, 4).Value = "Aporte"
}

$wsIngreso.Range("A473").Select()

# --- Sheet "Gastos": append a duplicate of row 48 (arbitro, agua y hielo) ---
$wsGastos = $wb.Worksheets.Item("Gastos")
$wsGastos.Cells.Item(49, 1).Value = 45101
$wsGastos.Cells.Item(49, 2).Value = "Arbitro, agua y hielo"
$wsGastos.Cells.Item(49, 3).Value = 940

$wsGastos.Range("A49").Select()

# --- Sheet "Cuentas por cobrar": append new debt row ---
$wsCuentas = $wb.Worksheets.Item("Cuentas por cobrar")
$wsCuentas.Cells.Item(3, 1).Copy()
$wsCuentas.Cells.Item(4, 1).PasteSpecial(-4122)
$wsCuentas.Cells.Item(4, 1).Value = 45109
$wsCuentas.Cells.Item(4, 2).Value = "Carlos"
$wsCuentas.Cells.Item(4, 3).Value = "Tecnica"
$wsCuentas.Cells.Item(4, 4).Value = 100
$wsCuentas.Cells.Item(4, 6).Value = "Le pitaron 3 segundos e hizo un pique"

$wsCuentas.Range("A5").Select()

$wsIngreso.Activate()
